$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: add a new sub-bullet "Is it possible to overfit?" right
# before the existing "How to tune the parameters?" bullet (paragraph
# 5 in the original document).
# ---------------------------------------------------------------------
$d.Paragraphs.Item(5).Range.InsertParagraphBefore() | Out-Null
$overfitPara = $d.Paragraphs.Item(5)
$overfitPara.Range.Text = "Is it possible to overfit?"
$overfitPara.Range.ListFormat.ListLevelNumber = 2

# ---------------------------------------------------------------------
# Hunk 2: after "Significant variables (and how these change with
# horizon)" -- turn the trailing empty bullet into a short "For
# tuning:" section followed by three new bullets, a blank paragraph
# and a closing paragraph of prose.
#
# After the hunk-1 insert above, the trailing empty bullet (originally
# paragraph 10) is now paragraph 11.
# ---------------------------------------------------------------------
$idx = 11

# Two blank (Normal-style) paragraphs before the trailing bullet.
$d.Paragraphs.Item($idx).Range.InsertParagraphBefore() | Out-Null
$idx = $idx + 1
$d.Paragraphs.Item($idx - 1).Style = "Normal"

$d.Paragraphs.Item($idx).Range.InsertParagraphBefore() | Out-Null
$idx = $idx + 1
$d.Paragraphs.Item($idx - 1).Style = "Normal"

# "For tuning:" paragraph before the trailing bullet.
$d.Paragraphs.Item($idx).Range.InsertParagraphBefore() | Out-Null
$idx = $idx + 1
$forTuning = $d.Paragraphs.Item($idx - 1)
$forTuning.Style = "Normal"
$forTuning.Range.Text = "For tuning:"

# $idx now points at the old trailing empty bullet; turn it into the
# first of three top-level (ilvl 0) bullets.
$trailing = $d.Paragraphs.Item($idx)
$trailing.Range.ListFormat.ListLevelNumber = 1
$trailing.Range.Text = "Pick a few hyperparameters"

$d.Paragraphs.Item($idx).Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$d.Paragraphs.Item($idx).Range.Text = "Select values for those hyperparameters"

$d.Paragraphs.Item($idx).Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$d.Paragraphs.Item($idx).Range.Text = "Use cross-validation to choose best values"

$d.Paragraphs.Item($idx).Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$d.Paragraphs.Item($idx).Style = "Normal"

$d.Paragraphs.Item($idx).Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$closing = $d.Paragraphs.Item($idx)
$closing.Style = "Normal"
$closing.Range.Text = "Random forest and cross-validation are inherently cross-sectional techniques. It doesn’t quite make sense to k-fold the time series data, because you’ll end up validating past data on future data. This seems problematic because the future data theoretically has no effect on the past data. For all we know, the time series is fundamentally different in the 2010s as opposed to the 1960s. Maybe there’s a structural break, maybe there’s some kind of non-linear time trend; at any rate, it would be best to avoid validating past predictions on future data."
